# Cosenza.xlsx update — add new days (2021-07-30 .. 2021-08-31) and correct
# a run of previously-provisional daily "Nuovi casi" figures (rows 413-459,
# row 494) with the final reconciled counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "Nuovi casi" — corrected historical daily counts (column C)
# ---------------------------------------------------------------------------
$wsNuoviCasi = $wb.Worksheets.Item("Nuovi casi")

$historicalFixes = @{
    413 = 247
    415 = 115
    416 = 156
    417 = 184
    418 = 107
    420 = 162
    421 = 85
    422 = 69
    423 = 99
    424 = 116
    425 = 188
    426 = 118
    428 = 47
    429 = 20
    432 = 122
    433 = 72
    435 = 11
    437 = 105
    439 = 72
    441 = 70
    445 = 61
    453 = 73
    494 = -2
}
foreach ($r in $historicalFixes.Keys) {
    $wsNuoviCasi.Cells.Item($r, 3).Value = $historicalFixes[$r]
}

# ---------------------------------------------------------------------------
# 2) New daily data rows 509-521 (2021-07-30 .. 2021-08-11) on every sheet,
#    then bare date rows 522-541 (2021-08-12 .. 2021-08-31, not yet reported)
# ---------------------------------------------------------------------------
$newValuesBySheet = @{
    "Nuovi casi"        = @{509=21;  510=35; 511=43;  512=1;  513=14;  514=62;  515=44;  516=53; 517=65;  518=22; 519=-19; 520=14; 521=81}
    "Deceduti"          = @{509=0;   510=0;  511=0;   512=1;  513=0;   514=2;   515=0;   516=0;  517=0;   518=0;  519=1;   520=0;  521=0}
    "Dimessi   Guariti" = @{509=39;  510=3;  511=26;  512=0;  513=137; 514=181; 515=329; 516=51; 517=75;  518=23; 519=15;  520=41; 521=56}
    "Ricoveri"          = @{509=32;  510=29; 511=30;  512=31; 513=31;  514=35;  515=36;  516=34; 517=32;  518=32; 519=33;  520=36; 521=39}
    "Terapia"           = @{509=4;   510=3;  511=2;   512=2;  513=2;   514=1;   515=1;   516=1;  517=1;   518=1;  519=1;   520=1;  521=1}
}

$lastDataRow = 521
$lastDateRow = 541
$firstNewDateSerial = 44409   # row 511's date; rows 509/510 already had dates

foreach ($sheetName in $wb.Worksheets) { }

foreach ($entry in $newValuesBySheet.GetEnumerator()) {
    $sheetName = $entry.Key
    $values = $entry.Value
    $ws = $wb.Worksheets.Item($sheetName)

    # Extend the date column first: rows 511-541 don't exist yet. Clone the
    # date-column formatting from row 510 (already formatted as dd/mm/yyyy).
    $ws.Range("A510").Copy()
    $ws.Range("A511:A541").PasteSpecial(-4122)
    for ($i = 0; $i -lt ($lastDateRow - 511 + 1); $i++) {
        $ws.Cells.Item(511 + $i, 1).Value = $firstNewDateSerial + $i
    }

    # Clone column C / D number formatting from row 508 down through the new
    # data rows (509-521) before writing values/formulas into them.
    $ws.Range("C508").Copy()
    $ws.Range("C509:C521").PasteSpecial(-4122)
    $ws.Range("D508").Copy()
    $ws.Range("D509:D521").PasteSpecial(-4122)

    foreach ($r in 509..$lastDataRow) {
        $ws.Cells.Item($r, 3).Value = $values[$r]
        $firstRow = $r - 6
        $ws.Cells.Item($r, 4).Formula = "=AVERAGE(C$firstRow`:C$r)"
    }
}

# ---------------------------------------------------------------------------
# 3) View state: active sheet moves from "Terapia" to "Deceduti"; each sheet
#    ends up scrolled/selected around the newly-entered block.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Nuovi casi").Range("A2:D521").Select()

$wsDeceduti = $wb.Worksheets.Item("Deceduti")
$wsDeceduti.Activate()
$excel.ActiveWindow.ScrollRow = 495
$wsDeceduti.Range("C509:C521").Select()

$wsGuariti = $wb.Worksheets.Item("Dimessi   Guariti")
$excel.ActiveWindow.ScrollRow = 501
$wsGuariti.Range("C509:C521").Select()

$wsRicoveri = $wb.Worksheets.Item("Ricoveri")
$excel.ActiveWindow.ScrollRow = 503
$wsRicoveri.Range("C509:C521").Select()

$wsTerapia = $wb.Worksheets.Item("Terapia")
$wsTerapia.Range("C509:C521").Select()

$wsDeceduti.Activate()
